# Add invalid currency fixture: new "Accrual currency" / "accrual exchange
# rate" columns (L, M) with sample rows for the first two data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 12).Value = "Accrual currency"
$ws.Cells.Item(1, 13).Value = "accrual exchange rate"

# Row 2 data
$ws.Cells.Item(2, 12).Value = "INRA"
$ws.Cells.Item(2, 13).Value = 1.223

# Row 3 data
$ws.Cells.Item(3, 12).Value = "A"
$ws.Cells.Item(3, 13).Value = 2.33

# Update selection to match the recorded edit state
$ws.Range("L1:M3").Select()
